$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue "D2" "96.395.30"
Set-TextValue "E2" "  +4.88%  "
Set-TextValue "D3" "3.599.95"
Set-TextValue "E3" "  +9.45%  "
Set-TextValue "E4" "  -0.01%  "
Set-TextValue "D5" "239.94"
Set-TextValue "E5" "  +5.56%  "
Set-TextValue "D6" "638.60"
Set-TextValue "E6" "  +4.83%  "
Set-TextValue "D7" "1.48"
Set-TextValue "E7" "  +8.98%  "
Set-TextValue "E8" "  +7.34%  "
Set-TextValue "E10" "  +9.78%  "
Set-TextValue "D11" "3.596.65"
Set-TextValue "E11" "  +9.51%  "
Set-TextValue "D12" "43.16"
Set-TextValue "E12" "  +4.49%  "
Set-TextValue "E13" "  +4.93%  "
Set-TextValue "D14" "6.42"
Set-TextValue "E14" "  +8.66%  "
Set-TextValue "D15" "4.284.65"
Set-TextValue "E15" "  +9.91%  "
Set-TextValue "D16" "96.338.84"
Set-TextValue "E16" "  +4.99%  "
Set-TextValue "E17" "  +6.14%  "
Set-TextValue "E18" "  +5.18%  "
Set-TextValue "D19" "3.606.27"
Set-TextValue "E19" "  +9.67%  "
Set-TextValue "D20" "13.21"
Set-TextValue "E20" "  +25.19%  "
Set-TextValue "D21" "18.10"
Set-TextValue "E21" "  +6.86%  "
Set-TextValue "D22" "0.499"
Set-TextValue "E22" "  +15.30%  "
Set-TextValue "D23" "516.64"
Set-TextValue "E23" "  +7.06%  "
Set-TextValue "E24" "  +3.33%  "
Set-TextValue "E25" "  +12.92%  "
Set-TextValue "D26" "6.67"
Set-TextValue "E26" "  +9.71%  "
Set-TextValue "D27" "97.12"
Set-TextValue "E28" "  +7.64%  "
Set-TextValue "E29" "  +19.14%  "
Set-TextValue "D30" "11.60"
Set-TextValue "E30" "  +7.33%  "
Set-TextValue "D31" "0.144"
Set-TextValue "E31" "  +5.95%  "
Set-TextValue "E32" "  -0.11%  "
Set-TextValue "E33" "  +7.46%  "
Set-TextValue "D34" "0.997"
Set-TextValue "E34" "  +0.12%  "
Set-TextValue "D35" "30.43"
Set-TextValue "E35" "  +10.13%  "
Set-TextValue "D36" "0.570"
Set-TextValue "E36" "  +11.29%  "
Set-TextValue "D37" "575.10"
Set-TextValue "E37" "  +7.34%  "
Set-TextValue "D38" "7.90"
Set-TextValue "E38" "  +9.31%  "
Set-TextValue "D39" "1.49"
Set-TextValue "E39" "  +11.50%  "
Set-TextValue "E40" "  +4.54%  "
Set-TextValue "E41" "  +0.02%  "
Set-TextValue "D42" "0.927"
Set-TextValue "E42" "  +9.59%  "
Set-TextValue "D43" "1.75"
Set-TextValue "E43" "  +6.74%  "
Set-TextValue "D44" "0.0431"
Set-TextValue "E44" "  +7.72%  "
Set-TextValue "E45" "  +0.39%  "
Set-TextValue "D46" "5.67"
Set-TextValue "E46" "  +8.30%  "
Set-TextValue "E47" "  +7.37%  "
Set-TextValue "E48" "  -0.36%  "
Set-TextValue "E49" "  +4.92%  "
Set-TextValue "D50" "8.12"
Set-TextValue "E50" "  +3.89%  "
Set-TextValue "E51" "  +6.97%  "
